# Update "城市市容环境卫生" (Urban Appearance & Environmental Sanitation) sheet:
# drop the oldest six years (2004年-2009年) and append the newest year (2021年),
# keeping the rolling 12-year window (2010年-2021年).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 2-7 (2004年..2009年); remaining rows shift up so 2010年 becomes row 2.
$ws.Rows("2:7").Delete()

# Copy the (now) last data row's formatting down onto the new row 13 so the
# year label A13 keeps the bordered/centered/bold style used by every other
# year cell in column A, then overwrite it with the new year + stats.
$ws.Range("A12").Copy($ws.Range("A13"))

$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 184063
$ws.Range("C13").Value = 327512
$ws.Range("D13").Value = 24869.205146
$ws.Range("F13").Value = 1034211.2
